$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats/styles) of the last existing data row
# down into the two new rows, so the new cells reuse the same style indices.
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16: 2025-10-08 (serial 45938), Sifangping station
$ws.Range("A16").Value = 45938
$ws.Range("B16").Value = "四方坪站"
$ws.Range("C16").Value = 11559.41
$ws.Range("D16").Value = 9475.15
$ws.Range("E16").Value = 4044.18
$ws.Range("F16").Value = 463

# Row 17: 2025-10-08 (serial 45938), Gaoling station
$ws.Range("A17").Value = 45938
$ws.Range("B17").Value = "高岭站"
$ws.Range("C17").Value = 5580.13
$ws.Range("D17").Value = 4610.29
$ws.Range("E17").Value = 1362.14
$ws.Range("F17").Value = 192

# Move the view down to the new rows and update the selection, matching
# the author's on-screen position after entering the new data.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
[void]$ws.Range("H18").Select()
